# Feria Lagunitas de Puerto Montt - Repollo: insert a new weekly record
# at row 706, pushing the existing rows 706:796 down to 707:797.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 706 (shifts rows 706..796 -> 707..797).
$ws.Rows.Item(706).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A706").Value = 4
$ws.Range("B706").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C706").Value = "Los Lagos"
$ws.Range("D706").Value = 45154
$ws.Range("E706").Value = 10
$ws.Range("F706").Value = 100112006
$ws.Range("G706").Value = "Repollo"
$ws.Range("H706").Value = "Crespo record"
$ws.Range("I706").Value = "Primera"
$ws.Range("J706").Value = 250
$ws.Range("K706").Value = 1500
$ws.Range("L706").Value = 1500
$ws.Range("M706").Value = 1500
$ws.Range("N706").Value = "`$/unidad"
$ws.Range("O706").Value = "Región Metropolitana"
$ws.Range("P706").Value = 1500
$ws.Range("Q706").Value = 1
$ws.Range("R706").Value = "Hortaliza"
